$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Rebecca / 9087 / goa
$ws.Range("A2").Value = "Rebecca"
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "9087"
$ws.Range("C2").Value = "goa"

# Row 3: kaviya / chennai / 123
$ws.Range("A3").Value = "kaviya"
$ws.Range("B3").Value = "chennai"
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "123"
